# Add GME and AMC to the companies list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Symbols first for both new rows, then fill in the rest (matches authoring order).
$ws.Range("A77").Value = "GME"
$ws.Range("A78").Value = "AMC"

$ws.Range("B78").Value = "AMC Entertainment Holdings Inc"
$ws.Range("C78").Value = "amc entertainment holdings, amc entertainment"

$ws.Range("B77").Value = "GameStop Corp"
$ws.Range("C77").Value = "gamestop, gamestop corporation"

# Match the style (left alignment) used by columns A and B in the rest of the sheet.
$ws.Range("A77:B78").HorizontalAlignment = -4131

# Update the selection / scroll position to mirror the authored workbook state.
$ws.Range("A78").Select()
$excel.ActiveWindow.ScrollRow = 57
